# Daten aktualisiert am 2024-03-07
# Append five new ticker rows right after the existing data in column A
# (existing data occupies A1:A363, so new rows go to A364:A368).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @("IMX-USD", "TAO-USD", "GRT-USD", "MNT-USD", "PEPE-USD")

$startRow = 364
for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value = $newTickers[$i]
}
